$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 3.9
$ws.Range("P5").Value = 1.71

$ws.Range("F7").Value = 1.92
$ws.Range("K7").Value = 3.65

$ws.Range("F8").Value = 1.43
$ws.Range("I8").Value = 8.6
$ws.Range("K8").Value = 6.2
$ws.Range("Q8").Value = 1.46

$ws.Range("F9").Value = 2.14
$ws.Range("G9").Value = 2.7
$ws.Range("K9").Value = 5.6

$ws.Range("T10").Value = 2.06

$ws.Range("H11").Value = 1.78

$ws.Range("P12").Value = 1.24

$ws.Range("P13").Value = 1.24

$ws.Range("G15").Value = 2.6
$ws.Range("H15").Value = 3.8
$ws.Range("K15").Value = 3.7
$ws.Range("P15").Value = 1.25

$ws.Range("P16").Value = 2.76
$ws.Range("Q16").Value = 1.31

$ws.Range("F17").Value = 3.7
$ws.Range("I17").Value = 2.02

$ws.Range("P18").Value = 1.24

$ws.Range("G19").Value = 2.7
$ws.Range("H19").Value = 3.25

$ws.Range("F20").Value = 2.6
$ws.Range("G20").Value = 2.96
$ws.Range("J20").Value = 3.15
$ws.Range("Q20").Value = 1.97

$ws.Range("G21").Value = 2.94
$ws.Range("P21").Value = 1.64
$ws.Range("T21").Value = 2.02
$ws.Range("AC21").Value = 6.8
$ws.Range("AH21").Value = 20

$ws.Range("F22").Value = 1.94
$ws.Range("G22").Value = 1.96
$ws.Range("P22").Value = 2.8
$ws.Range("Q22").Value = 1.53
$ws.Range("S22").Value = 2.32
$ws.Range("T22").Value = 1.54
$ws.Range("U22").Value = 2.8
$ws.Range("X22").Value = 26
$ws.Range("Y22").Value = 23
$ws.Range("AH22").Value = 13.5
$ws.Range("AI22").Value = 34
$ws.Range("AN22").Value = 8

$ws.Range("P23").Value = 1.76
$ws.Range("AA23").Value = 70
$ws.Range("AK23").Value = 29

$ws.Range("F24").Value = 2.82
$ws.Range("G24").Value = 3.05
$ws.Range("I24").Value = 2.86
$ws.Range("L24").Value = 1.41
$ws.Range("N24").Value = 3.7
$ws.Range("O24").Value = 1.32
$ws.Range("T24").Value = 1.59
$ws.Range("V24").Value = 1.53
$ws.Range("W24").Value = 1.5
$ws.Range("AH24").Value = 18.5

$ws.Range("G25").Value = 2.54
$ws.Range("M25").Value = 1.08
$ws.Range("S25").Value = 4.2
